$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "di"
$ws.Range("C2").Value = "d"
$ws.Range("E2").Value = 34
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 619.91

# Row 3 updates
$ws.Range("B3").Value = "dd"
$ws.Range("C3").Value = "rrt"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 51.48

# Row 4 update
$ws.Range("A4").Value = 4
